{"js": "// The \"Form Code\" column of the SOP forms table uses the prefix\n// \"F-SW-CR/xx\". This renames the series to \"F-SW-SD/xx\", and while\n// renumbering, rows that used to read 07/08 become 09/10.\n//\n// A couple of the old codes (\u2026/03 and \u2026/04) are split across two runs\n// in the original document (e.g. \"F-SW-CR/0\" + \"3\"), so rather than\n// trying to patch individual runs, we find the paragraph that holds\n// each code and overwrite its whole text \u2014 that also naturally\n// collapses the paragraph back down to a single run, matching the\n// target edit.\n\nconst mapping = [\n  [\"F-SW-CR/01\", \"F-SW-SD/01\"],\n  [\"F-SW-CR/02\", \"F-SW-SD/02\"],\n  [\"F-SW-CR/03\", \"F-SW-SD/03\"],\n  [\"F-SW-CR/04\", \"F-SW-SD/04\"],\n  [\"F-SW-CR/05\", \"F-SW-SD/05\"],\n  [\"F-SW-CR/06\", \"F-SW-SD/06\"],\n  [\"F-SW-CR/07\", \"F-SW-SD/09\"],\n  [\"F-SW-CR/08\", \"F-SW-SD/10\"],\n];\n\nconst body = context.document.body;\n\n// \"F-SW-CR/0\" is common to every old code, including the ones split\n// across two runs, so searching for that prefix reliably locates every\n// row regardless of how the text is chunked into runs.\nconst results = body.search(\"F-SW-CR/0\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst paragraphs = [];\nfor (let i = 0; i < results.items.length; i++) {\n  paragraphs.push(results.items[i].paragraphs.getFirst());\n}\nfor (let i = 0; i < paragraphs.length; i++) {\n  paragraphs[i].load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.length; i++) {\n  const oldText = paragraphs[i].text.trim();\n  const hit = mapping.find(([oldCode]) => oldCode === oldText);\n  if (hit) {\n    paragraphs[i].insertText(hit[1], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# The \"Form Code\" column of the SOP forms table uses the prefix\n# \"F-SW-CR/xx\". Rename the series to \"F-SW-SD/xx\"; while renumbering,\n# the rows that used to read 07/08 become 09/10.\n#\n# Word's Find/Replace operates on the logical paragraph text, so it\n# transparently handles the couple of old codes (\u2026/03 and \u2026/04) that\n# happen to be split across two runs in the original document - the\n# replacement collapses them back down to a single run.\n\n$d = $word.ActiveDocument\n\n$mapping = @(\n    @{ Old = \"F-SW-CR/01\"; New = \"F-SW-SD/01\" },\n    @{ Old = \"F-SW-CR/02\"; New = \"F-SW-SD/02\" },\n    @{ Old = \"F-SW-CR/03\"; New = \"F-SW-SD/03\" },\n    @{ Old = \"F-SW-CR/04\"; New = \"F-SW-SD/04\" },\n    @{ Old = \"F-SW-CR/05\"; New = \"F-SW-SD/05\" },\n    @{ Old = \"F-SW-CR/06\"; New = \"F-SW-SD/06\" },\n    @{ Old = \"F-SW-CR/07\"; New = \"F-SW-SD/09\" },\n    @{ Old = \"F-SW-CR/08\"; New = \"F-SW-SD/10\" }\n)\n\nforeach ($pair in $mapping) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute(\n        $null,   # FindText (already set above)\n        $false,  # MatchCase\n        $false,  # MatchWholeWord\n        $false,  # MatchWildcards\n        $false,  # MatchSoundsLike\n        $false,  # MatchAllWordForms\n        $true,   # Forward\n        1,       # Wrap (wdFindContinue)\n        $false,  # Format\n        $null,   # ReplaceWith (already set above)\n        2        # Replace (wdReplaceAll)\n    )\n}\n"}
